# Colocando header nos gráficos
# Adds a header label in column A (row 1) for each data table, fixes
# accented Portuguese words that were previously written without
# diacritics, removes the bold/bordered header style from the row-label
# column (A2:A12, etc.) on the per-source tables, removes the "Teto" row
# from the "Emissoes Totais" sheet, and updates the "Custo Total" sheet
# with a proper year header and refreshed cost figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the same visual style as the existing bold header cells
# (e.g. B1) onto a newly introduced A1 header cell, then set its text.
# ---------------------------------------------------------------------
function Set-HeaderCell($ws, [string]$addr, [string]$text) {
    $ws.Range("B1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($addr).Value = $text
}

# =====================================================================
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# =====================================================================
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia"

    # Remove the bold/bordered style from the technology labels and fix
    # the missing Portuguese accents.
    $ws.Range("A2").Style = "Normal"

    $ws.Range("A3").Style = "Normal"
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").Style = "Normal"
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").Style = "Normal"

    $ws.Range("A6").Style = "Normal"
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").Style = "Normal"

    $ws.Range("A8").Style = "Normal"
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").Style = "Normal"

    $ws.Range("A10").Style = "Normal"

    $ws.Range("A11").Style = "Normal"
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").Style = "Normal"
}

# =====================================================================
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# =====================================================================
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período"

$ws5.Range("A2").Style = "Normal"
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").Style = "Normal"
$ws5.Range("A3").Value = "P.Crítico"

# The "Teto" row is no longer part of the table.
$ws5.Rows.Item(4).Delete()

# =====================================================================
# Sheet 6: "Custo Total (bilhões de R$)"
# =====================================================================
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "Tipo Expansão"

# Re-use the "2015" text header (same style/content as on the other
# sheets) instead of typing a bare number, so it stays a text value.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)  # xlPasteAll

$ws6.Range("A2").Style = "Normal"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 570

$ws6.Range("A3").Style = "Normal"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
